$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "as of" timestamp on every sheet (09:24 -> 11:02)
# ---------------------------------------------------------------------------
$wsMD410 = $wb.Worksheets.Item("MD410 Attendance")
$wsMD410.Range("A1").Value = "MD410 Registrees as of 23/04/2021 11:02"

$ws410E = $wb.Worksheets.Item("410E Attendance")
$ws410E.Range("A1").Value = "410E Registrees as of 23/04/2021 11:02"

$ws410W = $wb.Worksheets.Item("410W Attendance")
$ws410W.Range("A1").Value = "410W Registrees as of 23/04/2021 11:02"

$ws410EVoting = $wb.Worksheets.Item("410E Voting")
$ws410EVoting.Range("A1").Value = "410E Voting details as of 23/04/2021 11:02"

$ws410WVoting = $wb.Worksheets.Item("410W Voting")
$ws410WVoting.Range("A1").Value = "410W Voting details as of 23/04/2021 11:02"

# ---------------------------------------------------------------------------
# 2. Insert a new registree row into the "410E Attendance" sheet.
#    The new row is inserted above the current row 23 (Erasmus, Freddie),
#    pushing all the following rows (and the footer summary rows) down by
#    one, and updating the attendee count.
# ---------------------------------------------------------------------------
$ws = $ws410E

$ws.Rows.Item(23).Insert()

# Pick up the same look & feel (style + row height) as the surrounding data
# rows by copying the formatting of the row directly above.
$ws.Range("A22:E22").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 25

$ws.Range("A23").Value = "Elske"
$ws.Range("B23").Value = "Valerie"
$ws.Range("C23").Value = "Krugersdorp"
$ws.Range("D23").Value = "No"
$ws.Range("E23").Value = "No"

# The "Number of attendees" footer row (now shifted down to row 133) needs
# its count bumped from 129 to 130 to reflect the newly added registree.
$ws.Range("A133").Value = "Number of attendees: 130"
